# Add a new worksheet named "rho" as the last sheet in the workbook,
# containing theta / Jxx / Jyy / beta / gamma / trace_sq data.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "rho"

# Header row (bold, thin border, centered horizontally, top-aligned vertically)
$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace_sq")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
    @(0,  0.01740544426373977,  0.9825945467348716,  -0.003085342615102135, 0, 0.9658140314412309),
    @(5,  0.05033371666978453,  0.9496659873001346,  -0.1835881672030711,   0, 0.9718082007424992),
    @(10, 0.1399931911095196,   0.8600067895532182,  -0.3350744138643648,   0, 0.983759497287755),
    @(15, 0.2704532149268831,   0.7295475547190609,  -0.4433701993202432,   0, 0.9985388433513923),
    @(20, 0.4291608892232774,   0.5708390444083918,  -0.4833663521194739,   0, 0.9773223441825746),
    @(25, 0.6097607456828397,   0.3902392247587773,  -0.479829167203665,    0, 0.9845668789147497),
    @(30, 0.753628888887208,    0.2463711759298113,  -0.4102873272298778,   0, 0.9653266402652793),
    @(35, 0.8858979724518178,   0.114102346096456,   -0.3007003160677003,   0, 0.9786759231453867),
    @(40, 0.9602931292450565,   0.0397068449576482,  -0.1601799746889957,   0, 0.9750547761944881),
    @(45, 0.982273357351477,    0.01772661936939805,  0.01315406549478772,  0, 0.9655212404748923),
    @(50, 0.9519487132597725,   0.04805128938414006,  0.169831850366278,    0, 0.9662009938860826),
    @(55, 0.8620490949993135,   0.1379509111830151,   0.323692642435192,    0, 0.9717129496187137),
    @(60, 0.7350963442737286,   0.2649036552338432,   0.4190852155444523,   0, 0.9618054176967311),
    @(65, 0.5750090264247757,   0.4249909747538475,   0.4798413337025017,   0, 0.971748120150985),
    @(70, 0.402011643665746,    0.5979883629923279,   0.4737040360781434,   0, 0.9679944715105249),
    @(75, 0.2526874659623621,   0.7473125403944333,   0.4167847893150157,   0, 0.9697461096939854),
    @(80, 0.1244082937189699,   0.8755917157896242,   0.3091404604805468,   0, 0.9732739249177326),
    @(85, 0.04396000833269664,  0.9560399951553169,   0.1579034642402508,   0, 0.9658119627073334),
    @(90, 0.01568535574023989,  0.9843146263341712,  -0.01022692347884475,  0, 0.9693304939277613)
)

$row = 2
foreach ($record in $data) {
    for ($col = 1; $col -le $record.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $record[$col - 1]
    }
    $row++
}

